$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit re-orders the data rows 2-34 (rows 8, 13, 17, 29 are unchanged).
# Excel has no atomic "move many rows at once" primitive that safely handles
# permutation cycles, so stage every source row that needs to move into a
# scratch area (+1000 rows down, well outside the used range) first, then
# clear each destination row (so cells that become blank really go blank,
# since Copy onto a non-empty cell with an empty source leaves it untouched),
# then copy each staged row into its real destination, and finally wipe the
# scratch area so the sheet dimension/used-range is unaffected.

# --- Stage source rows ---
$ws.Range("A2:AY2").Copy($ws.Range("A1002:AY1002"))
$ws.Range("A3:AY3").Copy($ws.Range("A1003:AY1003"))
$ws.Range("A4:AY4").Copy($ws.Range("A1004:AY1004"))
$ws.Range("A5:AY5").Copy($ws.Range("A1005:AY1005"))
$ws.Range("A6:AY6").Copy($ws.Range("A1006:AY1006"))
$ws.Range("A7:AY7").Copy($ws.Range("A1007:AY1007"))
$ws.Range("A9:AY9").Copy($ws.Range("A1009:AY1009"))
$ws.Range("A10:AY10").Copy($ws.Range("A1010:AY1010"))
$ws.Range("A11:AY11").Copy($ws.Range("A1011:AY1011"))
$ws.Range("A12:AY12").Copy($ws.Range("A1012:AY1012"))
$ws.Range("A14:AY14").Copy($ws.Range("A1014:AY1014"))
$ws.Range("A15:AY15").Copy($ws.Range("A1015:AY1015"))
$ws.Range("A16:AY16").Copy($ws.Range("A1016:AY1016"))
$ws.Range("A18:AY18").Copy($ws.Range("A1018:AY1018"))
$ws.Range("A19:AY19").Copy($ws.Range("A1019:AY1019"))
$ws.Range("A20:AY20").Copy($ws.Range("A1020:AY1020"))
$ws.Range("A21:AY21").Copy($ws.Range("A1021:AY1021"))
$ws.Range("A22:AY22").Copy($ws.Range("A1022:AY1022"))
$ws.Range("A23:AY23").Copy($ws.Range("A1023:AY1023"))
$ws.Range("A24:AY24").Copy($ws.Range("A1024:AY1024"))
$ws.Range("A25:AY25").Copy($ws.Range("A1025:AY1025"))
$ws.Range("A26:AY26").Copy($ws.Range("A1026:AY1026"))
$ws.Range("A27:AY27").Copy($ws.Range("A1027:AY1027"))
$ws.Range("A28:AY28").Copy($ws.Range("A1028:AY1028"))
$ws.Range("A30:AY30").Copy($ws.Range("A1030:AY1030"))
$ws.Range("A31:AY31").Copy($ws.Range("A1031:AY1031"))
$ws.Range("A32:AY32").Copy($ws.Range("A1032:AY1032"))
$ws.Range("A33:AY33").Copy($ws.Range("A1033:AY1033"))
$ws.Range("A34:AY34").Copy($ws.Range("A1034:AY1034"))

# --- Clear destination rows so blanked-out cells truly become blank ---
$ws.Range("A2:AY2").Clear()
$ws.Range("A3:AY3").Clear()
$ws.Range("A4:AY4").Clear()
$ws.Range("A5:AY5").Clear()
$ws.Range("A6:AY6").Clear()
$ws.Range("A7:AY7").Clear()
$ws.Range("A9:AY9").Clear()
$ws.Range("A10:AY10").Clear()
$ws.Range("A11:AY11").Clear()
$ws.Range("A12:AY12").Clear()
$ws.Range("A14:AY14").Clear()
$ws.Range("A15:AY15").Clear()
$ws.Range("A16:AY16").Clear()
$ws.Range("A18:AY18").Clear()
$ws.Range("A19:AY19").Clear()
$ws.Range("A20:AY20").Clear()
$ws.Range("A21:AY21").Clear()
$ws.Range("A22:AY22").Clear()
$ws.Range("A23:AY23").Clear()
$ws.Range("A24:AY24").Clear()
$ws.Range("A25:AY25").Clear()
$ws.Range("A26:AY26").Clear()
$ws.Range("A27:AY27").Clear()
$ws.Range("A28:AY28").Clear()
$ws.Range("A30:AY30").Clear()
$ws.Range("A31:AY31").Clear()
$ws.Range("A32:AY32").Clear()
$ws.Range("A33:AY33").Clear()
$ws.Range("A34:AY34").Clear()

# --- Copy staged rows into their final destination rows ---
$ws.Range("A1012:AY1012").Copy($ws.Range("A2:AY2"))
$ws.Range("A1011:AY1011").Copy($ws.Range("A3:AY3"))
$ws.Range("A1010:AY1010").Copy($ws.Range("A4:AY4"))
$ws.Range("A1002:AY1002").Copy($ws.Range("A5:AY5"))
$ws.Range("A1016:AY1016").Copy($ws.Range("A6:AY6"))
$ws.Range("A1015:AY1015").Copy($ws.Range("A7:AY7"))
$ws.Range("A1014:AY1014").Copy($ws.Range("A9:AY9"))
$ws.Range("A1004:AY1004").Copy($ws.Range("A10:AY10"))
$ws.Range("A1003:AY1003").Copy($ws.Range("A11:AY11"))
$ws.Range("A1007:AY1007").Copy($ws.Range("A12:AY12"))
$ws.Range("A1006:AY1006").Copy($ws.Range("A14:AY14"))
$ws.Range("A1005:AY1005").Copy($ws.Range("A15:AY15"))
$ws.Range("A1009:AY1009").Copy($ws.Range("A16:AY16"))
$ws.Range("A1034:AY1034").Copy($ws.Range("A18:AY18"))
$ws.Range("A1027:AY1027").Copy($ws.Range("A19:AY19"))
$ws.Range("A1031:AY1031").Copy($ws.Range("A20:AY20"))
$ws.Range("A1028:AY1028").Copy($ws.Range("A21:AY21"))
$ws.Range("A1021:AY1021").Copy($ws.Range("A22:AY22"))
$ws.Range("A1032:AY1032").Copy($ws.Range("A23:AY23"))
$ws.Range("A1025:AY1025").Copy($ws.Range("A24:AY24"))
$ws.Range("A1033:AY1033").Copy($ws.Range("A25:AY25"))
$ws.Range("A1020:AY1020").Copy($ws.Range("A26:AY26"))
$ws.Range("A1023:AY1023").Copy($ws.Range("A27:AY27"))
$ws.Range("A1024:AY1024").Copy($ws.Range("A28:AY28"))
$ws.Range("A1026:AY1026").Copy($ws.Range("A30:AY30"))
$ws.Range("A1022:AY1022").Copy($ws.Range("A31:AY31"))
$ws.Range("A1030:AY1030").Copy($ws.Range("A32:AY32"))
$ws.Range("A1019:AY1019").Copy($ws.Range("A33:AY33"))
$ws.Range("A1018:AY1018").Copy($ws.Range("A34:AY34"))

# --- Clean up the scratch/staging rows used for the permutation ---
$ws.Range("A1002:AY1002").Clear()
$ws.Range("A1003:AY1003").Clear()
$ws.Range("A1004:AY1004").Clear()
$ws.Range("A1005:AY1005").Clear()
$ws.Range("A1006:AY1006").Clear()
$ws.Range("A1007:AY1007").Clear()
$ws.Range("A1009:AY1009").Clear()
$ws.Range("A1010:AY1010").Clear()
$ws.Range("A1011:AY1011").Clear()
$ws.Range("A1012:AY1012").Clear()
$ws.Range("A1014:AY1014").Clear()
$ws.Range("A1015:AY1015").Clear()
$ws.Range("A1016:AY1016").Clear()
$ws.Range("A1018:AY1018").Clear()
$ws.Range("A1019:AY1019").Clear()
$ws.Range("A1020:AY1020").Clear()
$ws.Range("A1021:AY1021").Clear()
$ws.Range("A1022:AY1022").Clear()
$ws.Range("A1023:AY1023").Clear()
$ws.Range("A1024:AY1024").Clear()
$ws.Range("A1025:AY1025").Clear()
$ws.Range("A1026:AY1026").Clear()
$ws.Range("A1027:AY1027").Clear()
$ws.Range("A1028:AY1028").Clear()
$ws.Range("A1030:AY1030").Clear()
$ws.Range("A1031:AY1031").Clear()
$ws.Range("A1032:AY1032").Clear()
$ws.Range("A1033:AY1033").Clear()
$ws.Range("A1034:AY1034").Clear()
